# Fixes duration and assignment of sub-jobs in FullJobShopProblem
#
# The "production orders" sheet had a single example order row describing
# product_id 4 (Cola 2L). This updates it to describe product_id 13
# (Fanta 500ml) instead, with a different amount and derived calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("production orders")

# Make sure this is the active sheet (it was already the tab-selected one).
$ws.Activate()

# production_order_nr: P1 -> P4
$ws.Range("A2").Value2 = "P4"

# days_till_delivery: -1 -> 4
$ws.Range("B2").Value2 = 4

# product_id: 4 (Cola 2L) -> 13 (Fanta 500ml)
$ws.Range("C2").Value2 = 13

# amount: 8000 -> 7500
$ws.Range("D2").Value2 = 7500

# comments: update note to reflect the new product/amount
$ws.Range("E2").Value2 = "product_id 13 = Fanta 500ml"

# liters required: was =D2*2 (2 liters per bottle for Cola 2L),
# now =D2*0.5 (0.5 liters per bottle for Fanta 500ml)
$ws.Range("F2").Formula = "=D2*0.5"

# Bottling time in minutes (G2) keeps its original formula =D2*60/1000;
# it recalculates automatically from the new amount in D2.

# Update the selected/active cell on this sheet.
$ws.Range("I12").Select()
